# Daily attendance processing - 2025-12-02 07:04:49
# For every row in column G ("Recorded By") whose value is a comma-separated
# list that begins with "System", move that leading "System" entry to the
# end of the list (e.g. "System, foo@bar.com" -> "foo@bar.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }

    if ($val.StartsWith("System, ")) {
        $parts = $val -split ", "
        $newParts = $parts[1..($parts.Length - 1)] + $parts[0]
        $newVal = $newParts -join ", "
        $cell.Value = $newVal
    }
}
